$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "No significant differences detected between overview_home_page_20240721-155448.png and overview_home_page_20240803-143929.png."
$ws.Range("B2").Value = "Success"
$ws.Range("C2").Value = "2024-08-03 14:39:56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.22%"
$ws.Range("E2").Style = $ws.Range("D2").Style
$ws.Range("F2").Value = "overview_home_page_20240721-155448.png"
$ws.Range("G2").Value = "overview_home_page_20240803-143929.png"

# Update row 3
$ws.Range("A3").Value = "No significant differences detected between login_home_page_20240721-155451.png and login_home_page_20240803-143934.png."
$ws.Range("B3").Value = "Success"
$ws.Range("C3").Value = "2024-08-03 14:39:56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.22%"
$ws.Range("E3").Style = $ws.Range("D3").Style
$ws.Range("F3").Value = "login_home_page_20240721-155451.png"
$ws.Range("G3").Value = "login_home_page_20240803-143934.png"

# Delete rows 4 and 5
$ws.Range("A4:G5").EntireRow.Delete()
